# Applies the cryptos.xlsx data refresh described in the commit diff.
# Numeric-looking text values are prefixed with a literal single quote so
# Excel stores them as text (matching the source inlineStr cells) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.635.97"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "2.345.05"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'515.87"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'133.42"
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "2.342.50"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  +6.05%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  +5.78%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.758.46"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'23.61"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "56.686.23"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "2.345.18"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").Value = "'10.36"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "'318.56"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'60.60"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  +4.96%  "
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").Value = "'170.64"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +9.22%  "
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").Value = "'6.20"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "'18.18"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "'0.942"
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("E39").Value = "  +7.08%  "
$ws.Range("D40").Value = "'37.46"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").Value = "'0.378"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "'137.72"
$ws.Range("E42").Value = "  +8.72%  "
$ws.Range("D43").Value = "'3.55"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").Value = "'274.95"
$ws.Range("E44").Value = "  +9.93%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "'0.0927"
$ws.Range("E46").Value = "  +2.71%  "
$ws.Range("D47").Value = "'0.0501"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'0.557"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E49").Value = "  +4.04%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'16.69"
